# Update countries & provincias Spain
# Applies the data refresh captured in the commit: reorders a handful of
# countries (whose case counts now put them in a different rank) and
# refreshes the "last updated" timestamp plus several countries' stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 14:52"

# --- simple stat refreshes (no re-ranking) ---
# Row 70: Uzbekistan
$ws.Range("D70").Value = 225
$ws.Range("E70").Value = 1265

# Row 74: Lituania
$ws.Range("E74").Value = 1021
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 35

# --- Sri Lanka jumps ahead of Vietnam/Kenia/Guatemala (rows 116-119) ---
$ws.Range("A116").Value = "Sri Lanka"
$ws.Range("B116").Value = 269
$ws.Range("C116").Value = 15
$ws.Range("D116").Value = 91
$ws.Range("E116").Value = 171
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 7

$ws.Range("A117").Value = "Vietnam"
$ws.Range("B117").Value = 268
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 203
$ws.Range("E117").Value = 65
$ws.Range("F117").Value = 8
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0

$ws.Range("A118").Value = "Kenia"
$ws.Range("B118").Value = 262
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 60
$ws.Range("E118").Value = 190
$ws.Range("F118").Value = 2
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 12

$ws.Range("A119").Value = "Guatemala"
$ws.Range("B119").Value = 257
$ws.Range("C119").Value = 22
$ws.Range("D119").Value = 21
$ws.Range("E119").Value = 229
$ws.Range("F119").Value = 3
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 7

# --- Mozambique jumps ahead of Siria/San Martin (Parte Francesa) (rows 166-168) ---
$ws.Range("A166").Value = "Mozambique"
$ws.Range("B166").Value = 39
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 8
$ws.Range("E166").Value = 31
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

$ws.Range("A167").Value = "Siria"
$ws.Range("B167").Value = 38
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 5
$ws.Range("E167").Value = 31
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 2

$ws.Range("A168").Value = "San Martin (Parte Francesa)"
$ws.Range("B168").Value = 37
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 19
$ws.Range("E168").Value = 16
$ws.Range("F168").Value = 5
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 2

# --- Bonaire jumps ahead of Burundi/Butan; Sudan del Sur/Islas Virgenes shift
#     down past Santo Tome y Principe (which stays put) (rows 208-213) ---
$ws.Range("A208").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B208").Value = 5
$ws.Range("C208").Value = 2
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 5
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Burundi"
$ws.Range("B209").Value = 5
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 4
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Butan"
$ws.Range("B210").Value = 5
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 2
$ws.Range("E210").Value = 3
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211 (Santo Tome y Principe) is unchanged.

$ws.Range("A212").Value = "Sudan del Sur"
$ws.Range("B212").Value = 4
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 4
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 4
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 2
$ws.Range("E213").Value = 1
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 1
$ws.Range("H213").Value = 1
